# Add a new diary entry (row 36) to the "Impi" learning-diary workbook.
#
# Column layout (row 1 headers): A=PVM (date), B=Kello (time),
# C=Oppimisen sisältö (content), D=Oppimisen laatu (quality),
# E=Huomiot koodista (code notes), F=META, G=Tunnit (hours),
# H=Kertymä (running total formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the new row's cells in the same order the original author's Excel
# session produced them (A, B, D, C, F) so newly-appended shared-string
# entries land in the expected order.
$ws.Range("A36").Value = "27 marras"

$ws.Range("B36").Value = "18.00-20.45,21.15-22.00"
$ws.Range("B36").WrapText = $true
$ws.Range("B36").NumberFormat = "h:mm"

$ws.Range("D36").Value = "Vähän copy pasteksi loppupuolella meni. Pikkuhiljaa palastellaan eteenpäin samalla kun tehdään demoa, halusin käydä teorian nyt hieman ripeämmin jotta ehtii saada demon aikaiseksi kurssilla ja saa koodin kääntymään. Oppivelkaa kasvatettu kyllä, mutta eipä se suurimmalta osin mahdottomaltakaan tuntunut."
$ws.Range("D36").WrapText = $true

$ws.Range("C36").Value = "Törmäysfysiikkaa, törmäysten aiheuttamat voimat. 301-350"
$ws.Range("C36").WrapText = $true

$ws.Range("F36").Value = "Kirjasta kahlattu oleellisimmat asiat, nyt lueskellaan rauhallisempaan tahtiin relevantteja asioita hieman mitä ehtii. Kurssi alkaa olla loppusuoralla oppisisältötavoitteellisesti, nyt vikat demot pakettiin ja pientä koodikannan siisitimistä ja artikkelin etsintää."
$ws.Range("F36").WrapText = $true

# Hours logged for this entry; H3's SUM(G3:G60) formula recalculates
# automatically from 104.5 to 108.
$ws.Range("G36").Value = 3.5

# Match the row height Excel's wrap-text autofit produced for this entry.
$ws.Rows(36).RowHeight = 130.5

# Update the visible window/selection state to the end of the new entry.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B36").Select()
